$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and 1h volume-change figures to match the
# latest scrape. Values are stored as text (matching the source data feed),
# so each cell is marked as Text before the literal value is written; this
# prevents Excel from re-interpreting strings such as "-4.74%" or "307.95"
# as numbers/percentages.
$updates = @{
    'D2' = '307.95'
    'E2' = '-4.74%'
    'D3' = '49.33'
    'E3' = '-0.40%'
    'D4' = '5.173'
    'E4' = '-3.16%'
    'D5' = '0.07719'
    'E5' = '-5.46%'
    'D6' = '4.516'
    'E6' = '-1.81%'
    'E7' = '13.27%'
    'D8' = '1.549'
    'E8' = '-7.59%'
    'D9' = '0.1231'
    'E9' = '-8.56%'
    'D10' = '0.1946'
    'E10' = '-1.22%'
    'D11' = '0.09404'
    'E11' = '-3.51%'
    'D12' = '0.04646'
    'E12' = '5.18%'
    'D13' = '0.1047'
    'E13' = '-0.11%'
    'D14' = '0.001260'
    'E14' = '-5.59%'
    'D15' = '0.04178'
    'E15' = '-2.95%'
    'D16' = '0.005860'
    'E16' = '-0.59%'
    'E17' = '2,023.02%'
    'D18' = '3.334'
    'E18' = '-1.46%'
    'D19' = '2.233'
    'E19' = '-8.38%'
    'D20' = '0.3486'
    'E20' = '2.70%'
    'D21' = '7.963'
    'E21' = '-2.16%'
    'D22' = '0.1343'
    'E22' = '-5.33%'
    'D23' = '0.3038'
    'E23' = '4.41%'
    'E24' = '-2.44%'
    'D25' = '0.004005'
    'E25' = '-6.22%'
    'D26' = '0.0001354'
    'E26' = '0.36%'
    'D38' = '0.02582'
    'E38' = '-6.68%'
    'D39' = '0.05823'
    'E39' = '4.13%'
    'D40' = '0.01076'
    'E40' = '70.85%'
    'D41' = '0.007934'
    'E41' = '3.27%'
    'D42' = '0.1419'
    'E42' = '-2.05%'
    'D43' = '0.008466'
    'E43' = '10.32%'
    'D44' = '0.007698'
    'E44' = '-4.93%'
    'D45' = '0.3373'
    'E45' = '-4.21%'
    'D46' = '0.00007025'
    'E46' = '0.88%'
    'E47' = '0.36%'
    'D48' = '0.04882'
    'E48' = '-20.39%'
    'E49' = '0.21%'
    'D50' = '0.00002106'
    'E50' = '0.36%'
    'D51' = '0.0002006'
    'E51' = '0.36%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
